$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2999.889
$ws.Range("I113").Value = 1999.6666
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1999.6666
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 1254.3334
$ws.Range("N113").Value = -10008
$ws.Range("H116").Value = 33336400
$ws.Range("I116").Value = 200000000
$ws.Range("J116").Value = 3680
$ws.Range("K116").Value = 200000000
$ws.Range("L116").Value = 3680
$ws.Range("M116").Value = -199996558
$ws.Range("N116").Value = -10564
$ws.Range("H132").Value = 1729.0962
$ws.Range("I132").Value = 1249.425
$ws.Range("J132").Value = 3328
$ws.Range("K132").Value = 3748.275
$ws.Range("L132").Value = 9984
$ws.Range("M132").Value = -1218.275
$ws.Range("N132").Value = -15044
$ws.Range("H137").Value = 1839.4
$ws.Range("I137").Value = 1582.6111
$ws.Range("J137").Value = 2224.5833
$ws.Range("K137").Value = 4747.8333
$ws.Range("L137").Value = 6673.749899999999
$ws.Range("M137").Value = -2197.8333
$ws.Range("N137").Value = -11773.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12651.968
$ws.Range("I32").Value = 12694.742
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 12694.742
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -12407.742
$ws.Range("N32").Value = -10574
$ws.Range("H63").Value = 3475.8635
$ws.Range("I63").Value = 3897.6155
$ws.Range("J63").Value = 2866.6667
$ws.Range("K63").Value = 3897.6155
$ws.Range("L63").Value = 2866.6667
$ws.Range("M63").Value = -3211.6155
$ws.Range("N63").Value = -4238.6667
$ws.Range("H66").Value = 3475.8635
$ws.Range("I66").Value = 3897.6155
$ws.Range("J66").Value = 2866.6667
$ws.Range("K66").Value = 19488.0775
$ws.Range("L66").Value = 14333.3335
$ws.Range("M66").Value = -16056.0775
$ws.Range("N66").Value = -21197.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1316.5714
$ws.Range("I94").Value = 1049
$ws.Range("J94").Value = 1673.3334
$ws.Range("K94").Value = 1049
$ws.Range("L94").Value = 1673.3334
$ws.Range("M94").Value = -598
$ws.Range("N94").Value = -2575.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8250
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 9285.714
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 9285.714
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -9509.714
$ws.Range("H31").Value = 23259230
$ws.Range("I31").Value = 50003040
$ws.Range("J31").Value = 3742.5652
$ws.Range("K31").Value = 50003040
$ws.Range("L31").Value = 3742.5652
$ws.Range("M31").Value = -50002745
$ws.Range("N31").Value = -4332.5652
$ws.Range("H34").Value = 23259230
$ws.Range("I34").Value = 50003040
$ws.Range("J34").Value = 3742.5652
$ws.Range("K34").Value = 50003040
$ws.Range("L34").Value = 3742.5652
$ws.Range("M34").Value = -50002838
$ws.Range("N34").Value = -4146.5652
$ws.Range("H58").Value = 1284.9143
$ws.Range("I58").Value = 1199.2903
$ws.Range("J58").Value = 1948.5
$ws.Range("K58").Value = 1199.2903
$ws.Range("L58").Value = 1948.5
$ws.Range("M58").Value = -996.2902999999999
$ws.Range("N58").Value = -2354.5
$ws.Range("H74").Value = 31438
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 31438
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 31438
$ws.Range("N74").Value = -33186
$ws.Range("H77").Value = 31438
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 31438
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 94314
$ws.Range("N77").Value = -103050
$ws.Range("H104").Value = 16095.333
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 16095.333
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 16095.333
$ws.Range("N104").Value = -21337.333
$ws.Range("H122").Value = 1414.5555
$ws.Range("I122").Value = 1297.0834
$ws.Range("J122").Value = 1649.5
$ws.Range("K122").Value = 3891.2502
$ws.Range("L122").Value = 4948.5
$ws.Range("M122").Value = -1441.2502
$ws.Range("N122").Value = -9848.5
$ws.Range("H132").Value = 1634.3334
$ws.Range("I132").Value = 1336.3125
$ws.Range("J132").Value = 2588
$ws.Range("K132").Value = 4008.9375
$ws.Range("L132").Value = 7764
$ws.Range("M132").Value = -1478.9375
$ws.Range("N132").Value = -12824
$ws.Range("H136").Value = 1284.9143
$ws.Range("I136").Value = 1199.2903
$ws.Range("J136").Value = 1948.5
$ws.Range("K136").Value = 3597.8709
$ws.Range("L136").Value = 5845.5
$ws.Range("M136").Value = -1047.8709
$ws.Range("N136").Value = -10945.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 798.5
$ws.Range("I69").Value = 347.5
$ws.Range("J69").Value = 1249.5
$ws.Range("K69").Value = 1042.5
$ws.Range("L69").Value = 3748.5
$ws.Range("M69").Value = -231.5
$ws.Range("N69").Value = -5370.5
$ws.Range("H70").Value = 12980
$ws.Range("I70").Value = 15225
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 45675
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -45360
$ws.Range("N70").Value = -12630
$ws.Range("H72").Value = 798.5
$ws.Range("I72").Value = 347.5
$ws.Range("J72").Value = 1249.5
$ws.Range("K72").Value = 3127.5
$ws.Range("L72").Value = 11245.5
$ws.Range("M72").Value = 928.5
$ws.Range("N72").Value = -19357.5
$ws.Range("H73").Value = 12980
$ws.Range("I73").Value = 15225
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 45675
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -44583
$ws.Range("N73").Value = -14184
$ws.Range("H74").Value = 3250
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -14122
$ws.Range("H77").Value = 3250
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 36000
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -46608

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 26166.166
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 26166.166
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 26166.166
$ws.Range("N97").Value = -28148.166
$ws.Range("H101").Value = 182720.67
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 182720.67
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 182720.67
$ws.Range("N101").Value = -189210.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22730482
$ws.Range("I122").Value = 31253290
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 93759870
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -93757420
$ws.Range("N122").Value = -13885

